$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "OrderDetailsData"
Write-Output $ws1.Name
